$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A93").Value = "chile"
$ws.Range("C93").Select() | Out-Null
